# "Schedule FIFO working perfectly"
#
# 1) Parameters!B8 (number block sections B): 2 -> 4
# 2) TimetableComplete (sheet2): the FIFO schedule advances one step.
#    The previously-sorted queue (L6:Q11, 6 rows) loses its earliest
#    entry (old L7:Q7, value 3 in its first column) which goes back to
#    the small "next up" holding slot (now J5:O5). The remaining 5
#    rows become the new primary table A2:F6 (sorted ascending by the
#    3rd column), replacing the old primary table (old A2:G3, 2 rows).

$wb = $excel.ActiveWorkbook

# --- Sheet "Parameters" ---
$params = $wb.Worksheets.Item("Parameters")
$params.Range("B8").Value = 4

# --- Sheet "TimetableComplete" ---
$ws = $wb.Worksheets.Item("TimetableComplete")

# Wipe the old primary table and the old sorted queue (contents AND
# formatting) so rows that end up fully blank disappear entirely,
# same as in the target layout.
$ws.Range("A2:G3").Clear()
$ws.Range("L6:Q11").Clear()

# New primary table, A2:F3 -- bold (first two rows, "locked in").
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "IC"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 42
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 52

$ws.Range("A3").Value = 13
$ws.Range("B3").Value = "IC"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 54
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 8

$ws.Range("A2:F3").Font.Bold = $true
$ws.Range("G2").Font.Bold = $true
$ws.Range("G3").Font.Bold = $true

# Rows 4-6, not bold.
$ws.Range("A4").Value = 13
$ws.Range("B4").Value = "IC"
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 10

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "IC"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 52

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "IC"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 12

# The entry that fell out of the sorted table goes back to the small
# "next up" holding slot, now at J5:O5 (was L7:Q7) -- bold, as before.
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = "IC"
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 48
$ws.Range("N5").Value = 6
$ws.Range("O5").Value = 3
$ws.Range("J5:P5").Font.Bold = $true

# Sheet/view bookkeeping to match the refreshed layout.
$ws.Range("E20").Select()

$wb.Save()
